$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40,8).Value = 12492.6875  # H40
$ws.Cells.Item(40,9).Value = 10000  # I40
$ws.Cells.Item(40,10).Value = 12658.866  # J40
$ws.Cells.Item(40,11).Value = 10000  # K40
$ws.Cells.Item(40,12).Value = 12658.866  # L40
$ws.Cells.Item(40,13).Value = -9825  # M40
$ws.Cells.Item(40,14).Value = -13008.866  # N40
$ws.Cells.Item(41,8).Value = 1457.5294  # H41
$ws.Cells.Item(41,9).Value = 1705.3334  # I41
$ws.Cells.Item(41,11).Value = 1705.3334  # K41
$ws.Cells.Item(41,13).Value = -1265.3334  # M41
$ws.Cells.Item(46,8).Value = 23999.666  # H46
$ws.Cells.Item(46,9).Value = 11499.5  # I46
$ws.Cells.Item(46,11).Value = 34498.5  # K46
$ws.Cells.Item(46,13).Value = -34379.5  # M46
$ws.Cells.Item(60,8).Value = 23999.666  # H60
$ws.Cells.Item(60,9).Value = 11499.5  # I60
$ws.Cells.Item(60,11).Value = 34498.5  # K60
$ws.Cells.Item(60,13).Value = -34014.5  # M60
$ws.Cells.Item(110,8).Value = 67991.664  # H110
$ws.Cells.Item(110,10).Value = 67991.664  # J110
$ws.Cells.Item(110,12).Value = 67991.664  # L110
$ws.Cells.Item(110,14).Value = -76171.664  # N110
$ws.Cells.Item(133,8).Value = 88138.78  # H133
$ws.Cells.Item(133,10).Value = 88138.78  # J133
$ws.Cells.Item(133,12).Value = 88138.78  # L133
$ws.Cells.Item(133,14).Value = -98258.78  # N133
$ws.Cells.Item(134,8).Value = 53105.555  # H134
$ws.Cells.Item(134,10).Value = 55743.75  # J134
$ws.Cells.Item(134,12).Value = 55743.75  # L134
$ws.Cells.Item(134,14).Value = -65883.75  # N134
$ws.Cells.Item(136,8).Value = 96491.664  # H136
$ws.Cells.Item(136,10).Value = 96491.664  # J136
$ws.Cells.Item(136,12).Value = 96491.664  # L136
$ws.Cells.Item(136,14).Value = -106691.664  # N136

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4,8).Value = 454.625  # H4
$ws.Cells.Item(4,9).Value = 460.5  # I4
$ws.Cells.Item(4,10).Value = 437  # J4
$ws.Cells.Item(4,11).Value = 460.5  # K4
$ws.Cells.Item(4,12).Value = 437  # L4
$ws.Cells.Item(4,13).Value = -344.5  # M4
$ws.Cells.Item(4,14).Value = -669  # N4
$ws.Cells.Item(45,8).Value = 2791.4443  # H45
$ws.Cells.Item(45,9).Value = 2446.1428  # I45
$ws.Cells.Item(45,11).Value = 2446.1428  # K45
$ws.Cells.Item(45,13).Value = -2069.1428  # M45
$ws.Cells.Item(102,8).Value = 69217.766  # H102
$ws.Cells.Item(102,9).Value = 73102.92999999999  # I102
$ws.Cells.Item(102,10).Value = 51087  # J102
$ws.Cells.Item(102,11).Value = 73102.92999999999  # K102
$ws.Cells.Item(102,12).Value = 51087  # L102
$ws.Cells.Item(102,13).Value = -71480.92999999999  # M102
$ws.Cells.Item(102,14).Value = -54331  # N102
$ws.Cells.Item(132,8).Value = 2140.0645  # H132
$ws.Cells.Item(132,9).Value = 1621.2106  # I132
$ws.Cells.Item(132,10).Value = 2961.5833  # J132
$ws.Cells.Item(132,11).Value = 4863.6318  # K132
$ws.Cells.Item(132,12).Value = 8884.749899999999  # L132
$ws.Cells.Item(132,13).Value = -2333.6318  # M132
$ws.Cells.Item(132,14).Value = -13944.7499  # N132
$ws.Cells.Item(138,8).Value = 83747.25  # H138
$ws.Cells.Item(138,10).Value = 83747.25  # J138
$ws.Cells.Item(138,12).Value = 83747.25  # L138
$ws.Cells.Item(138,14).Value = -94027.25  # N138

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(5,8).Value = 679.6667  # H5
$ws.Cells.Item(5,9).Value = 759.6667  # I5
$ws.Cells.Item(5,11).Value = 759.6667  # K5
$ws.Cells.Item(5,13).Value = -646.6667  # M5
$ws.Cells.Item(22,8).Value = 65385.938  # H22
$ws.Cells.Item(22,9).Value = 91444.63  # I22
$ws.Cells.Item(22,10).Value = 8056.8  # J22
$ws.Cells.Item(22,11).Value = 91444.63  # K22
$ws.Cells.Item(22,12).Value = 8056.8  # L22
$ws.Cells.Item(22,13).Value = -91271.63  # M22
$ws.Cells.Item(22,14).Value = -8402.799999999999  # N22
$ws.Cells.Item(50,8).Value = 77020  # H50
$ws.Cells.Item(50,10).Value = 77020  # J50
$ws.Cells.Item(50,12).Value = 77020  # L50
$ws.Cells.Item(50,14).Value = -78168  # N50
$ws.Cells.Item(51,8).Value = 72809.836  # H51
$ws.Cells.Item(51,10).Value = 72809.836  # J51
$ws.Cells.Item(51,12).Value = 72809.836  # L51
$ws.Cells.Item(51,14).Value = -73791.836  # N51
$ws.Cells.Item(52,8).Value = 61251.11  # H52
$ws.Cells.Item(52,10).Value = 61251.11  # J52
$ws.Cells.Item(52,12).Value = 61251.11  # L52
$ws.Cells.Item(52,14).Value = -61777.11  # N52
$ws.Cells.Item(53,8).Value = 34093.8  # H53
$ws.Cells.Item(53,10).Value = 32940  # J53
$ws.Cells.Item(53,12).Value = 32940  # L53
$ws.Cells.Item(53,14).Value = -34088  # N53
$ws.Cells.Item(55,8).Value = 34334  # H55
$ws.Cells.Item(55,10).Value = 34334  # J55
$ws.Cells.Item(55,12).Value = 34334  # L55
$ws.Cells.Item(55,14).Value = -34880  # N55
$ws.Cells.Item(109,8).Value = 99990  # H109
$ws.Cells.Item(109,10).Value = 99990  # J109
$ws.Cells.Item(109,12).Value = 99990  # L109
$ws.Cells.Item(109,14).Value = -102764  # N109
$ws.Cells.Item(119,8).Value = 49726.285  # H119
$ws.Cells.Item(119,10).Value = 49726.285  # J119
$ws.Cells.Item(119,12).Value = 49726.285  # L119
$ws.Cells.Item(119,14).Value = -59402.285  # N119
$ws.Cells.Item(121,8).Value = 61251.11  # H121
$ws.Cells.Item(121,10).Value = 61251.11  # J121
$ws.Cells.Item(121,12).Value = 61251.11  # L121
$ws.Cells.Item(121,14).Value = -64745.11  # N121
$ws.Cells.Item(134,8).Value = 1559.3721  # H134
$ws.Cells.Item(134,9).Value = 1185.2162  # I134
$ws.Cells.Item(134,11).Value = 3555.6486  # K134
$ws.Cells.Item(134,13).Value = -1020.6486  # M134
$ws.Cells.Item(135,8).Value = 43272.184  # H135
$ws.Cells.Item(135,10).Value = 43272.184  # J135
$ws.Cells.Item(135,12).Value = 43272.184  # L135
$ws.Cells.Item(135,14).Value = -53412.184  # N135
$ws.Cells.Item(138,8).Value = 91720.75  # H138
$ws.Cells.Item(138,10).Value = 91720.75  # J138
$ws.Cells.Item(138,12).Value = 91720.75  # L138
$ws.Cells.Item(138,14).Value = -102000.75  # N138
$ws.Cells.Item(140,8).Value = 99990  # H140
$ws.Cells.Item(140,10).Value = 99990  # J140
$ws.Cells.Item(140,12).Value = 99990  # L140
$ws.Cells.Item(140,14).Value = -110350  # N140

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(9,8).Value = 37973.11  # H9
$ws.Cells.Item(9,10).Value = 37973.11  # J9
$ws.Cells.Item(9,12).Value = 37973.11  # L9
$ws.Cells.Item(9,14).Value = -38309.11  # N9
$ws.Cells.Item(16,8).Value = 2736.0715  # H16
$ws.Cells.Item(16,9).Value = 2596  # I16
$ws.Cells.Item(16,10).Value = 3249.6667  # J16
$ws.Cells.Item(16,11).Value = 2596  # K16
$ws.Cells.Item(16,12).Value = 3249.6667  # L16
$ws.Cells.Item(16,13).Value = -2309  # M16
$ws.Cells.Item(16,14).Value = -3823.6667  # N16
$ws.Cells.Item(18,8).Value = 24462.715  # H18
$ws.Cells.Item(18,10).Value = 23936.334  # J18
$ws.Cells.Item(18,12).Value = 23936.334  # L18
$ws.Cells.Item(18,14).Value = -24396.334  # N18
$ws.Cells.Item(113,8).Value = 2736.0715  # H113
$ws.Cells.Item(113,9).Value = 2596  # I113
$ws.Cells.Item(113,10).Value = 3249.6667  # J113
$ws.Cells.Item(113,11).Value = 2596  # K113
$ws.Cells.Item(113,12).Value = 3249.6667  # L113
$ws.Cells.Item(113,13).Value = -426  # M113
$ws.Cells.Item(113,14).Value = -7589.6667  # N113
$ws.Cells.Item(138,8).Value = 80050.875  # H138
$ws.Cells.Item(138,10).Value = 80050.875  # J138
$ws.Cells.Item(138,12).Value = 80050.875  # L138
$ws.Cells.Item(138,14).Value = -90330.875  # N138

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8,8).Value = 207090.8  # H8
$ws.Cells.Item(8,9).Value = 207090.8  # I8
$ws.Cells.Item(8,11).Value = 621272.3999999999  # K8
$ws.Cells.Item(8,13).Value = -621133.3999999999  # M8
$ws.Cells.Item(122,8).Value = 1263106.8  # H122
$ws.Cells.Item(122,9).Value = 722  # I122
$ws.Cells.Item(122,11).Value = 6498  # K122
$ws.Cells.Item(122,13).Value = -4048  # M122

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102,8).Value = 979.4286  # H102
$ws.Cells.Item(102,9).Value = 742.8333  # I102
$ws.Cells.Item(102,10).Value = 2399  # J102
$ws.Cells.Item(102,11).Value = 742.8333  # K102
$ws.Cells.Item(102,12).Value = 2399  # L102
$ws.Cells.Item(102,13).Value = 879.1667  # M102
$ws.Cells.Item(102,14).Value = -5643  # N102
$ws.Cells.Item(122,8).Value = 20057.934  # H122
$ws.Cells.Item(122,9).Value = 37325  # I122
$ws.Cells.Item(122,10).Value = 15741.167  # J122
$ws.Cells.Item(122,11).Value = 111975  # K122
$ws.Cells.Item(122,12).Value = 47223.501  # L122
$ws.Cells.Item(122,13).Value = -109525  # M122
$ws.Cells.Item(122,14).Value = -52123.501  # N122
$ws.Cells.Item(123,8).Value = 38997.6  # H123
$ws.Cells.Item(123,10).Value = 38997.6  # J123
$ws.Cells.Item(123,12).Value = 38997.6  # L123
$ws.Cells.Item(123,14).Value = -43897.6  # N123
$ws.Cells.Item(126,8).Value = 3553.2778  # H126
$ws.Cells.Item(126,9).Value = 2953.625  # I126
$ws.Cells.Item(126,11).Value = 8860.875  # K126
$ws.Cells.Item(126,13).Value = -6390.875  # M126
$ws.Cells.Item(132,8).Value = 7083.3335  # H132
$ws.Cells.Item(132,9).Value = 7300  # I132
$ws.Cells.Item(132,10).Value = 6812.5  # J132
$ws.Cells.Item(132,11).Value = 21900  # K132
$ws.Cells.Item(132,12).Value = 20437.5  # L132
$ws.Cells.Item(132,13).Value = -19370  # M132
$ws.Cells.Item(132,14).Value = -25497.5  # N132
$ws.Cells.Item(135,8).Value = 95121.664  # H135
$ws.Cells.Item(135,10).Value = 95121.664  # J135
$ws.Cells.Item(135,12).Value = 95121.664  # L135
$ws.Cells.Item(135,14).Value = -105261.664  # N135

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16,8).Value = 3283.2856  # H16
$ws.Cells.Item(16,9).Value = 3283.2856  # I16
$ws.Cells.Item(16,10).Value = 0  # J16
$ws.Cells.Item(16,11).Value = 3283.2856  # K16
$ws.Cells.Item(16,12).Value = 0  # L16
$ws.Cells.Item(16,13).Value = -3113.2856  # M16
$ws.Cells.Item(16,14).ClearContents()  # N16
$ws.Cells.Item(55,8).Value = 1746.3096  # H55
$ws.Cells.Item(55,9).Value = 931.89655  # I55
$ws.Cells.Item(55,11).Value = 931.89655  # K55
$ws.Cells.Item(55,13).Value = -758.89655  # M55
$ws.Cells.Item(100,8).Value = 75600.60000000001  # H100
$ws.Cells.Item(100,9).Value = 93749.75  # I100
$ws.Cells.Item(100,11).Value = 93749.75  # K100
$ws.Cells.Item(100,13).Value = -93208.75  # M100
$ws.Cells.Item(108,8).Value = 78685.25  # H108
$ws.Cells.Item(108,10).Value = 78685.25  # J108
$ws.Cells.Item(108,12).Value = 78685.25  # L108
$ws.Cells.Item(108,14).Value = -86365.25  # N108
$ws.Cells.Item(116,8).Value = 266333.34  # H116
$ws.Cells.Item(116,10).Value = 266333.34  # J116
$ws.Cells.Item(116,12).Value = 266333.34  # L116
$ws.Cells.Item(116,14).Value = -275511.34  # N116
$ws.Cells.Item(123,8).Value = 81493.336  # H123
$ws.Cells.Item(123,10).Value = 81493.336  # J123
$ws.Cells.Item(123,12).Value = 81493.336  # L123
$ws.Cells.Item(123,14).Value = -91293.336  # N123

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(118,8).Value = 74800  # H118
$ws.Cells.Item(118,10).Value = 74800  # J118
$ws.Cells.Item(118,12).Value = 74800  # L118
$ws.Cells.Item(118,14).Value = -78114  # N118
$ws.Cells.Item(119,8).Value = 69998  # H119
$ws.Cells.Item(119,10).Value = 69998  # J119
$ws.Cells.Item(119,12).Value = 69998  # L119
$ws.Cells.Item(119,14).Value = -79674  # N119
$ws.Cells.Item(120,8).Value = 79800  # H120
$ws.Cells.Item(120,10).Value = 79800  # J120
$ws.Cells.Item(120,12).Value = 79800  # L120
$ws.Cells.Item(120,14).Value = -89476  # N120
$ws.Cells.Item(126,8).Value = 1590.6857  # H126
$ws.Cells.Item(126,9).Value = 1286.8636  # I126
$ws.Cells.Item(126,11).Value = 3860.5908  # K126
$ws.Cells.Item(126,13).Value = -1390.5908  # M126
